$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4423792296040858
$ws.Range("C2").Value = 0.04231202262400302
$ws.Range("D2").Value = 0.1787246208188265
$ws.Range("E2").Value = 0.1640199045123225
$ws.Range("F2").Value = 1.578264622813151
$ws.Range("J2").Value = 0.1910775826062974
$ws.Range("K2").Value = 0.4076749453711557
$ws.Range("O2").Value = 3.9257739601521
$ws.Range("B3").Value = 0.4028768011925195
$ws.Range("C3").Value = 0.03704680571156871
$ws.Range("D3").Value = 0.1726011840276698
$ws.Range("E3").Value = 0.1597661717657815
$ws.Range("F3").Value = 1.582086128929888
$ws.Range("J3").Value = 0.1871676593355787
$ws.Range("K3").Value = 0.3655392386035885
$ws.Range("O3").Value = 3.949525278875143
$ws.Range("B4").Value = 0.3787075689793085
$ws.Range("C4").Value = 0.03379939988658975
$ws.Range("D4").Value = 0.1689111441245785
$ws.Range("E4").Value = 0.1572353702693228
$ws.Range("F4").Value = 1.58528031314134
$ws.Range("J4").Value = 0.1848750716527832
$ws.Range("K4").Value = 0.3396974074448451
$ws.Range("O4").Value = 3.966251166974217
$ws.Range("B5").Value = 0.368880430009483
$ws.Range("C5").Value = 0.03247247588018354
$ws.Range("D5").Value = 0.1674250754426936
$ws.Range("E5").Value = 0.1562244653391573
$ws.Range("F5").Value = 1.586795165354474
$ws.Range("J5").Value = 0.1839680236684913
$ws.Range("K5").Value = 0.3291746831738465
$ws.Range("O5").Value = 3.973605674330756
$ws.Range("B6").Value = 0.3672499869540218
$ws.Range("C6").Value = 0.03225192722202053
$ws.Range("D6").Value = 0.1671793841977518
$ws.Range("E6").Value = 0.1560578398908419
$ws.Range("F6").Value = 1.587059582515401
$ws.Range("J6").Value = 0.1838190522284648
$ws.Range("K6").Value = 0.3274278948557594
$ws.Range("O6").Value = 3.974859408084541
$ws.Range("B7").Value = 0.3785749468726465
$ws.Range("C7").Value = 0.03378151891638481
$ws.Range("D7").Value = 0.1688910308653675
$ws.Range("E7").Value = 0.1572216541202458
$ws.Range("F7").Value = 1.585299879721866
$ws.Range("J7").Value = 0.184862728737258
$ws.Range("K7").Value = 0.3395554610234512
$ws.Range("O7").Value = 3.966348172011095
$ws.Range("B8").Value = 0.4287413743881245
$ws.Range("C8").Value = 0.04049963822154723
$ws.Range("D8").Value = 0.1765988398988725
$ws.Range("E8").Value = 0.1625364309417492
$ws.Range("F8").Value = 1.579406332968901
$ws.Range("J8").Value = 0.1897070061627559
$ws.Range("K8").Value = 0.3931407266033204
$ws.Range("O8").Value = 3.933518635278062
$ws.Range("B9").Value = 0.5277767394654518
$ws.Range("C9").Value = 0.05355572955569698
$ws.Range("D9").Value = 0.1922636791973531
$ws.Range("E9").Value = 0.1736000934153026
$ws.Range("F9").Value = 1.574576491666562
$ws.Range("J9").Value = 0.2000647204693422
$ws.Range("K9").Value = 0.4984368148983549
$ws.Range("O9").Value = 3.886149601908329
$ws.Range("B10").Value = 0.60092241782732
$ws.Range("C10").Value = 0.06307327762745274
$ws.Range("D10").Value = 0.2041041074354268
$ws.Range("E10").Value = 0.1821188158427134
$ws.Range("F10").Value = 1.57513259589642
$ws.Range("J10").Value = 0.2081990403756038
$ws.Range("K10").Value = 0.575910493964841
$ws.Range("O10").Value = 3.861734003133478
$ws.Range("B11").Value = 0.6342785473484867
$ws.Range("C11").Value = 0.06738630401692092
$ws.Range("D11").Value = 0.2095618957753089
$ws.Range("E11").Value = 0.1860788460550751
$ws.Range("F11").Value = 1.576277746034719
$ws.Range("J11").Value = 0.2120138081150316
$ws.Range("K11").Value = 0.6111762820209492
$ws.Range("O11").Value = 3.852885921631923
$ws.Range("B12").Value = 0.6469209721408902
$ws.Range("C12").Value = 0.06901709096939612
$ws.Range("D12").Value = 0.2116388137232406
$ws.Range("E12").Value = 0.187590571865762
$ws.Range("F12").Value = 1.576839710641721
$ws.Range("J12").Value = 0.2134748211788775
$ws.Range("K12").Value = 0.6245333086611424
$ws.Range("O12").Value = 3.849860459284002
$ws.Range("B13").Value = 0.6441977116143107
$ws.Range("C13").Value = 0.06866598264389268
$ws.Range("D13").Value = 0.2111910616671793
$ws.Range("E13").Value = 0.1872644551427953
$ws.Range("F13").Value = 1.576712973964078
$ws.Range("J13").Value = 0.2131594347294055
$ws.Range("K13").Value = 0.6216565264165865
$ws.Range("O13").Value = 3.850497582169453
$ws.Range("B14").Value = 0.6353184264411027
$ws.Range("C14").Value = 0.06752051966215333
$ws.Range("D14").Value = 0.2097325618451293
$ws.Range("E14").Value = 0.1862029734568011
$ws.Range("F14").Value = 1.576321407300881
$ws.Range("J14").Value = 0.2121336770534725
$ws.Range("K14").Value = 0.6122751232827
$ws.Range("O14").Value = 3.852630497445659
$ws.Range("B15").Value = 0.6298810470248952
$ws.Range("C15").Value = 0.06681856739191971
$ws.Range("D15").Value = 0.2088405099519548
$ws.Range("E15").Value = 0.1855543658938927
$ws.Range("F15").Value = 1.576098273409187
$ws.Range("J15").Value = 0.2115075119383079
$ws.Range("K15").Value = 0.6065290683348508
$ws.Range("O15").Value = 3.853979318458386
$ws.Range("B16").Value = 0.5987441126297881
$ws.Range("C16").Value = 0.06279107107090454
$ws.Range("D16").Value = 0.2037488576576578
$ws.Range("E16").Value = 0.1818617211947213
$ws.Range("F16").Value = 1.575075712846569
$ws.Range("J16").Value = 0.2079520379316477
$ws.Range("K16").Value = 0.5736061943568131
$ws.Range("O16").Value = 3.862357723876926
$ws.Range("B17").Value = 0.5796631481692032
$ws.Range("C17").Value = 0.06031603404332486
$ws.Range("D17").Value = 0.2006435314555404
$ws.Range("E17").Value = 0.1796180931909319
$ws.Range("F17").Value = 1.574676932284262
$ws.Range("J17").Value = 0.2058001690944309
$ws.Range("K17").Value = 0.5534144331505217
$ws.Range("O17").Value = 3.868076335990111
$ws.Range("B18").Value = 0.5686960153501275
$ws.Range("C18").Value = 0.05889090493838012
$ws.Range("D18").Value = 0.1988641682830234
$ws.Range("E18").Value = 0.1783356043504654
$ws.Range("F18").Value = 1.574531538433803
$ws.Range("J18").Value = 0.204573241924237
$ws.Range("K18").Value = 0.5418028320916903
$ws.Range("O18").Value = 3.871578110790551
$ws.Range("B19").Value = 0.5649840785100366
$ws.Range("C19").Value = 0.05840811644573307
$ws.Range("D19").Value = 0.1982628672242441
$ws.Range("E19").Value = 0.1779027488643194
$ws.Range("F19").Value = 1.574496732381164
$ws.Range("J19").Value = 0.2041596756995148
$ws.Range("K19").Value = 0.5378717342828168
$ws.Range("O19").Value = 3.872800252813335
$ws.Range("B20").Value = 0.5816935532509433
$ws.Range("C20").Value = 0.0605796673150536
$ws.Range("D20").Value = 0.2009734021541902
$ws.Range("E20").Value = 0.1798561049879623
$ws.Range("F20").Value = 1.574710691940723
$ws.Range("J20").Value = 0.2060281245398414
$ws.Range("K20").Value = 0.5555636609011003
$ws.Range("O20").Value = 3.867445576396051
$ws.Range("B21").Value = 0.6379261888495762
$ws.Range("C21").Value = 0.06785703739033977
$ws.Range("D21").Value = 0.2101606834274463
$ws.Range("E21").Value = 0.1865144270463261
$ws.Range("F21").Value = 1.576432937108848
$ws.Range("J21").Value = 0.2124345207683547
$ws.Range("K21").Value = 0.6150306020811911
$ws.Range("O21").Value = 3.851995182507693
$ws.Range("B22").Value = 0.6747423798062187
$ws.Range("C22").Value = 0.07259882822218344
$ws.Range("D22").Value = 0.2162243408207587
$ws.Range("E22").Value = 0.1909368123763286
$ws.Range("F22").Value = 1.578306489725179
$ws.Range("J22").Value = 0.2167173100586695
$ws.Range("K22").Value = 0.6539107793159928
$ws.Range("O22").Value = 3.843792616769122
$ws.Range("B23").Value = 0.6550871370451716
$ws.Range("C23").Value = 0.07006938897283987
$ws.Range("D23").Value = 0.2129826695103674
$ws.Range("E23").Value = 0.1885700412286226
$ws.Range("F23").Value = 1.577238091639529
$ws.Range("J23").Value = 0.2144227391765554
$ws.Range("K23").Value = 0.6331585183455957
$ws.Range("O23").Value = 3.847996968989037
$ws.Range("B24").Value = 0.5807755984493497
$ws.Range("C24").Value = 0.06046048555121786
$ws.Range("D24").Value = 0.2008242491326229
$ws.Range("E24").Value = 0.1797484767929589
$ws.Range("F24").Value = 1.574695167948448
$ws.Range("J24").Value = 0.2059250340771683
$ws.Range("K24").Value = 0.5545920046108677
$ws.Range("O24").Value = 3.867730076011384
$ws.Range("B25").Value = 0.5009162959379978
$ws.Range("C25").Value = 0.05003664357511184
$ws.Range("D25").Value = 0.1879674467499939
$ws.Range("E25").Value = 0.1705384921364796
$ws.Range("F25").Value = 1.57516252640464
$ws.Range("J25").Value = 0.1971706747311117
$ws.Range("K25").Value = 0.4699303097600875
$ws.Range("O25").Value = 3.89714143565709
